$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new data (ДОБЫЧА УРАНОВОЙ И ТОРИЕВОЙ РУД) ---
$ws.Range("A7").Value = "ДОБЫЧА УРАНОВОЙ И ТОРИЕВОЙ РУД"
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 139.155311205598
$ws.Range("D7").Value = 55.733366288235402
$ws.Range("E7").Value = 169.90504433138199
$ws.Range("F7").Value = 45691
$ws.Range("G7").Value = 4.7782654746083697
$ws.Range("H7").Value = 81.575167989218698
$ws.Range("I7").Value = 45659

# --- Row 8: new data (ДОБЫЧА МЕТАЛЛИЧЕСКИХ РУД) ---
$ws.Range("A8").Value = "ДОБЫЧА МЕТАЛЛИЧЕСКИХ РУД"
$ws.Range("B8").Value = 13
$ws.Range("C8").Value = 155.66508405208899
$ws.Range("D8").Value = 88.7403960742175
$ws.Range("E8").Value = 146.18927291052199
$ws.Range("F8").Value = 72.336215535268096
$ws.Range("G8").Value = 120.76426693053899
$ws.Range("H8").Value = 109.23168598919401
$ws.Range("I8").Value = 86.529342249587998

# --- Row 9: new data (ДОБЫЧА ПРОЧИХ ПОЛЕЗНЫХ ИСКОПАЕМЫХ) ---
$ws.Range("A9").Value = "ДОБЫЧА ПРОЧИХ ПОЛЕЗНЫХ ИСКОПАЕМЫХ"
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 87.342761340226602
$ws.Range("D9").Value = 166.24011583657099
$ws.Range("E9").Value = 190.01301006619099
$ws.Range("F9").Value = 72.009509943385197
$ws.Range("G9").Value = 131.62958829441601
$ws.Range("H9").Value = 105.239136610673
$ws.Range("I9").Value = 78.062207283610505
$ws.Rows.Item(9).RowHeight = 25.5

# --- Row 10: new data (ПРОИЗВОДСТВО ПИЩЕВЫХ ПРОДУКТОВ, ВКЛЮЧАЯ НАПИТКИ) ---
$ws.Range("A10").Value = "ПРОИЗВОДСТВО ПИЩЕВЫХ ПРОДУКТОВ, ВКЛЮЧАЯ НАПИТКИ"
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 102.48784689007
$ws.Range("D10").Value = 103.679188341931
$ws.Range("E10").Value = 88.905993801519898
$ws.Range("F10").Value = 120.928731280402
$ws.Range("G10").Value = 92.894728262206499
$ws.Range("H10").Value = 78.307898597556999
$ws.Range("I10").Value = 110.002230673667
$ws.Rows.Item(10).RowHeight = 25.5

# --- Row 6: fix punctuation in the shared string (semicolon -> comma) ---
# Done last so this edited text becomes the newest/trailing unique string,
# matching the target shared-strings ordering.
$ws.Range("A6").Value = "ДОБЫЧА СЫРОЙ НЕФТИ И ПРИРОДНОГО ГАЗА, ПРЕДОСТАВЛЕНИЕ УСЛУГ В  ЭТИХ ОБЛАСТЯХ"

# --- Selection moves to A15 ---
$ws.Range("A15").Select()
